# chamadaBelaVista.xlsx - "mudança sistema de filtro sem painel lateral"
#
# 1) Alunos!13: fill in the previously-missing Idade/Categoria/Data de
#    Nascimento for "carla camuratti" (I13, J13, M13).
# 2) Registros: two new "chamada" date columns (Y, Z) are appended, and
#    three new students (rows 6-8) are added to the roster with their
#    attendance marks in those two new columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Alunos sheet
# ---------------------------------------------------------------------
$alunos = $wb.Worksheets.Item("Alunos")

$alunos.Range("I13").Value = 51
$alunos.Range("J13").Value = "Não definida"

# Match the existing "Data de Nascimento" formatting used by the other
# rows in column M before writing the date, so the new cell reuses the
# same number format instead of minting a brand new one.
$alunos.Range("M13").NumberFormat = $alunos.Range("M12").NumberFormat
$alunos.Range("M13").Value = "05/09/1974"

# ---------------------------------------------------------------------
# 2) Registros sheet
# ---------------------------------------------------------------------
$registros = $wb.Worksheets.Item("Registros")

# New "chamada" (attendance) date columns, styled like the existing
# header cells (bold / centered / bordered), appended after column X.
$registros.Range("X1").Copy()
$registros.Range("Y1").PasteSpecial(-4122)
$registros.Range("Z1").PasteSpecial(-4122)
$registros.Range("Y1").Value = "19/12/2025"
$registros.Range("Z1").Value = "24/12/2025"

# Existing students get blank attendance cells in the two new columns.
$registros.Range("Y2").Value = "x"
$registros.Range("Y2").Value = ""
$registros.Range("Z2").Value = "x"
$registros.Range("Z2").Value = ""

$registros.Range("Y3").Value = "x"
$registros.Range("Y3").Value = ""
$registros.Range("Z3").Value = "x"
$registros.Range("Z3").Value = ""

$registros.Range("Y4").Value = "x"
$registros.Range("Y4").Value = ""
$registros.Range("Z4").Value = "x"
$registros.Range("Z4").Value = ""

$registros.Range("Y5").Value = "x"
$registros.Range("Y5").Value = ""
$registros.Range("Z5").Value = "x"
$registros.Range("Z5").Value = ""

# Three new students appended to the roster, with attendance recorded
# only for the two newly added dates (columns Y and Z).
$registros.Range("B6").Value = "joão do pão"
$registros.Range("Y6").Value = "f"
$registros.Range("Z6").Value = "j"

$registros.Range("B7").Value = "fernando lando"
$registros.Range("Y7").Value = "c"
$registros.Range("Z7").Value = "c"

$registros.Range("B8").Value = "ana cintra"
$registros.Range("Y8").Value = "c"
$registros.Range("Z8").Value = "c"

Write-Output "edit complete"
